$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 13334604
$ws.Range("I33").Value = 17545184
$ws.Range("K33").Value = 17545184
$ws.Range("M33").Value = -17544955

# Row 131
$ws.Range("H131").Value = 1798.1364
$ws.Range("I131").Value = 1508
$ws.Range("K131").Value = 4524
$ws.Range("M131").Value = 516

# Row 137
$ws.Range("H137").Value = 2056276.6
$ws.Range("I137").Value = 6341825
$ws.Range("J137").Value = 6666.4346
$ws.Range("K137").Value = 19025475
$ws.Range("L137").Value = 19999.3038
$ws.Range("M137").Value = -19022925
$ws.Range("N137").Value = -25099.3038

# Row 138
$ws.Range("H138").Value = 9292.933999999999
$ws.Range("J138").Value = 9292.933999999999
$ws.Range("L138").Value = 27878.802
$ws.Range("N138").Value = -38158.802

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3336.3635
$ws.Range("I2").Value = 3200
$ws.Range("K2").Value = 3200
$ws.Range("M2").Value = -3087

# Row 32
$ws.Range("H32").Value = 22575.56
$ws.Range("I32").Value = 15866.1875
$ws.Range("K32").Value = 15866.1875
$ws.Range("M32").Value = -15579.1875

# Row 61
$ws.Range("H61").Value = 3911.342
$ws.Range("I61").Value = 2731.3667
$ws.Range("K61").Value = 2731.3667
$ws.Range("M61").Value = -2519.3667

# Row 97
$ws.Range("H97").Value = 807.3913
$ws.Range("I97").Value = 700.5
$ws.Range("K97").Value = 700.5
$ws.Range("M97").Value = -204.5

# Row 102
$ws.Range("H102").Value = 387848.66
$ws.Range("I102").Value = 437915.97
$ws.Range("K102").Value = 437915.97
$ws.Range("M102").Value = -436293.97

# Row 110
$ws.Range("H110").Value = 1135.875
$ws.Range("I110").Value = 1135.875
$ws.Range("K110").Value = 1135.875
$ws.Range("M110").Value = 909.125

# Row 116
$ws.Range("H116").Value = 3336.3635
$ws.Range("I116").Value = 3200
$ws.Range("K116").Value = 3200
$ws.Range("M116").Value = -906

# Row 122
$ws.Range("H122").Value = 6963
$ws.Range("I122").Value = 7695.923
$ws.Range("J122").Value = 2199
$ws.Range("K122").Value = 23087.769
$ws.Range("L122").Value = 6597
$ws.Range("M122").Value = -20637.769
$ws.Range("N122").Value = -11497

# Row 136
$ws.Range("H136").Value = 3911.342
$ws.Range("I136").Value = 2731.3667
$ws.Range("K136").Value = 8194.1001
$ws.Range("M136").Value = -5644.1001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3336.3635
$ws.Range("I3").Value = 3200
$ws.Range("K3").Value = 3200
$ws.Range("M3").Value = -3086

# Row 20
$ws.Range("H20").Value = 2488
$ws.Range("I20").Value = 2316.6667
$ws.Range("K20").Value = 2316.6667
$ws.Range("M20").Value = -2069.6667

# Row 99
$ws.Range("H99").Value = 2781
$ws.Range("I99").Value = 2322.5715
$ws.Range("K99").Value = 2322.5715
$ws.Range("M99").Value = -824.5715

# Row 107
$ws.Range("H107").Value = 16365.593
$ws.Range("I107").Value = 17008.24
$ws.Range("K107").Value = 17008.24
$ws.Range("M107").Value = -15088.24

# Row 134
$ws.Range("H134").Value = 3005.111
$ws.Range("I134").Value = 2043.625
$ws.Range("K134").Value = 6130.875
$ws.Range("M134").Value = -3595.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5105.102
$ws.Range("I31").Value = 3927.0833
$ws.Range("J31").Value = 6236
$ws.Range("K31").Value = 3927.0833
$ws.Range("L31").Value = 6236
$ws.Range("M31").Value = -3632.0833
$ws.Range("N31").Value = -6826

# Row 34
$ws.Range("H34").Value = 5105.102
$ws.Range("I34").Value = 3927.0833
$ws.Range("J34").Value = 6236
$ws.Range("K34").Value = 3927.0833
$ws.Range("L34").Value = 6236
$ws.Range("M34").Value = -3725.0833
$ws.Range("N34").Value = -6640

# Row 37
$ws.Range("H37").Value = 14730.63
$ws.Range("I37").Value = 13611.95
$ws.Range("J37").Value = 17926.857
$ws.Range("K37").Value = 13611.95
$ws.Range("L37").Value = 17926.857
$ws.Range("M37").Value = -13504.95
$ws.Range("N37").Value = -18140.857

# Row 134
$ws.Range("H134").Value = 2889.8
$ws.Range("I134").Value = 2572.9707
$ws.Range("J134").Value = 3237.2903
$ws.Range("K134").Value = 7718.9121
$ws.Range("L134").Value = 9711.8709
$ws.Range("M134").Value = -5183.9121
$ws.Range("N134").Value = -14781.8709

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 387.7647
$ws.Range("J2").Value = 831.5714
$ws.Range("L2").Value = 4989.428400000001
$ws.Range("N2").Value = -5215.428400000001

# Row 3
$ws.Range("H3").Value = 18365.5
$ws.Range("I3").Value = 15508
$ws.Range("J3").Value = 25033
$ws.Range("K3").Value = 46524
$ws.Range("L3").Value = 75099
$ws.Range("M3").Value = -46412
$ws.Range("N3").Value = -75323

# Row 34
$ws.Range("H34").Value = 2788.8333
$ws.Range("I34").Value = 46.666668
$ws.Range("J34").Value = 5531
$ws.Range("K34").Value = 140.000004
$ws.Range("L34").Value = 16593
$ws.Range("M34").Value = -56.00000399999999
$ws.Range("N34").Value = -16761

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 133
$ws.Range("H133").Value = 25000
$ws.Range("I133").Value = 20000
$ws.Range("K133").Value = 60000
$ws.Range("M133").Value = -54940

# Row 134
$ws.Range("H134").Value = 9737.056
$ws.Range("I134").Value = 10145.429
$ws.Range("K134").Value = 30436.287
$ws.Range("M134").Value = -25366.287

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 7599500.5
$ws.Range("I102").Value = 9461116
$ws.Range("J102").Value = 9837.923000000001
$ws.Range("K102").Value = 9461116
$ws.Range("L102").Value = 9837.923000000001
$ws.Range("M102").Value = -9459494
$ws.Range("N102").Value = -13081.923

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 913416.75
$ws.Range("I22").Value = 2195.8333
$ws.Range("K22").Value = 2195.8333
$ws.Range("M22").Value = -1900.8333

# Row 27
$ws.Range("H27").Value = 913416.75
$ws.Range("I27").Value = 2195.8333
$ws.Range("K27").Value = 2195.8333
$ws.Range("M27").Value = -2088.8333

# Row 61
$ws.Range("H61").Value = 6238.316
$ws.Range("I61").Value = 4291.385
$ws.Range("J61").Value = 10456.667
$ws.Range("K61").Value = 4291.385
$ws.Range("L61").Value = 10456.667
$ws.Range("M61").Value = -4089.385
$ws.Range("N61").Value = -10860.667

# Row 113
$ws.Range("H113").Value = 6238.316
$ws.Range("I113").Value = 4291.385
$ws.Range("J113").Value = 10456.667
$ws.Range("K113").Value = 4291.385
$ws.Range("L113").Value = 10456.667
$ws.Range("M113").Value = -2121.385
$ws.Range("N113").Value = -14796.667

# Row 132
$ws.Range("H132").Value = 4705.3716
$ws.Range("I132").Value = 4011
$ws.Range("J132").Value = 7482.857
$ws.Range("K132").Value = 12033
$ws.Range("L132").Value = 22448.571
$ws.Range("M132").Value = -9503
$ws.Range("N132").Value = -27508.571

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 83334190
$ws.Range("I107").Value = 111111944
$ws.Range("J107").Value = 55556436
$ws.Range("K107").Value = 333335832
$ws.Range("L107").Value = 166669308
$ws.Range("M107").Value = -333333912
$ws.Range("N107").Value = -166673148

# Row 132
$ws.Range("H132").Value = 7698.143
$ws.Range("I132").Value = 14901
$ws.Range("K132").Value = 44703
$ws.Range("M132").Value = -42173

# Row 136
$ws.Range("H136").Value = 14709628
$ws.Range("I136").Value = 16670530
$ws.Range("J136").Value = 2867.5
$ws.Range("K136").Value = 50011590
$ws.Range("L136").Value = 8602.5
$ws.Range("M136").Value = -50009040
$ws.Range("N136").Value = -13702.5

Write-Output "Done applying Phoenix_Profits updates"
